# Finish updating course objectives - close issue #5
#
# The "Statistical learning: cross-validation" topic is replaced with
# "Statistical learning: resampling methods", and the
# "Text analysis: topic modeling" topic is replaced with the merged
# "Text analysis: classification and topic modeling" (folding the old,
# now-removed standalone "topic modeling" entry into the classification
# week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "Statistical learning: resampling methods"
$ws.Range("D18").Value = "Text analysis: classification and topic modeling"

# Leave the selection where the author's last edit landed.
$ws.Range("D19").Select()
